$d = $word.ActiveDocument

$replacements = @(
    @('2024-01-24 Wednesday', '2024-01-25 Thursday'),
    @('33×20=660', '55×20=1100'),
    @('24×44=1056', '64×80=5120'),
    @('65×62=4030', '76×92=6992'),
    @('40×53=2120', '83×37=3071'),
    @('92×98=9016', '28×49=1372'),
    @('53×64=3392', '44×18=792'),
    @('38×93=3534', '63×92=5796'),
    @('84×79=6636', '18×57=1026'),
    @('99×39=3861', '15×71=1065'),
    @('75×64=4800', '37×49=1813'),
    @('75×96=7200', '25×26=650'),
    @('59×74=4366', '64×31=1984'),
    @('54×99=5346', '89×59=5251'),
    @('26×35=910', '76×58=4408'),
    @('54×41=2214', '22×31=682'),
    @('84×24=2016', '26×56=1456'),
    @('30×73=2190', '46×35=1610'),
    @('49×93=4557', '56×21=1176'),
    @('62×54=3348', '24×54=1296'),
    @('40×87=3480', '39×43=1677'),
    @('36×65=2340', '71×24=1704'),
    @('43×57=2451', '48×88=4224'),
    @('18×55=990', '45×51=2295'),
    @('18×93=1674', '59×85=5015'),
    @('28×20=560', '51×34=1734'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Done applying $($replacements.Count) replacements"
